# Commit: "Update column names from physical quantity to measurand"
#
# The header row relabels three columns:
#   F1: "Quantity"             -> "Measurand"
#   I1: "Physical quantity I"  -> "Measurand Level I"
#   J1: "Physical quantity II" -> "Measurand Level II"
#
# All other data cells keep their existing values/strings untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Measurand"
$ws.Range("I1").Value = "Measurand Level I"
$ws.Range("J1").Value = "Measurand Level II"

# Move the active selection to F1 (matches the author's cursor position
# after editing the "Measurand" header cell).
$ws.Range("F1").Select()
